$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '42.366.54'
$ws.Range('E2').Value = '  -0.47%  '
Set-TextCell 'D3' '2.180.43'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextCell 'D5' '252.04'
$ws.Range('E5').Value = '  +4.85%  '
Set-TextCell 'D6' '0.613'
$ws.Range('E6').Value = '  -1.07%  '
Set-TextCell 'D7' '73.18'
$ws.Range('E7').Value = '  -2.40%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -3.45%  '
Set-TextCell 'D10' '40.06'
$ws.Range('E10').Value = '  -3.28%  '
Set-TextCell 'D11' '0.0911'
$ws.Range('E11').Value = '  -1.83%  '
$ws.Range('E12').Value = '  -0.05%  '
Set-TextCell 'D13' '6.74'
$ws.Range('E13').Value = '  -2.17%  '
Set-TextCell 'D14' '2.505.68'
$ws.Range('E14').Value = '  -1.69%  '
Set-TextCell 'D15' '14.12'
$ws.Range('E15').Value = '  -3.59%  '
Set-TextCell 'D16' '2.177.46'
$ws.Range('E16').Value = '  -1.99%  '
Set-TextCell 'D17' '0.768'
$ws.Range('E17').Value = '  -4.37%  '
Set-TextCell 'D18' '42.283.43'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('E19').Value = '  -3.15%  '
Set-TextCell 'D20' '70.53'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -1.62%  '
Set-TextCell 'D22' '226.35'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('E23').Value = '  -6.09%  '
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -4.77%  '
Set-TextCell 'D27' '3.37'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D28' '2.22'
$ws.Range('E28').Value = '  +6.04%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D29' '2.16'
$ws.Range('E29').Value = '  -2.57%  '
Set-TextCell 'D30' '170.36'
$ws.Range('E30').Value = '  -1.20%  '
Set-TextCell 'D31' '36.63'
$ws.Range('E31').Value = '  +0.29%  '
Set-TextCell 'D32' '19.99'
$ws.Range('E32').Value = '  -1.41%  '
Set-TextCell 'D33' '0.0812'
$ws.Range('E33').Value = '  +1.98%  '
Set-TextCell 'D34' '5.08'
$ws.Range('E34').Value = '  -5.73%  '
Set-TextCell 'D35' '0.120'
$ws.Range('E35').Value = '  -1.73%  '
Set-TextCell 'D36' '0.107'
$ws.Range('E36').Value = '  -1.23%  '
Set-TextCell 'D37' '4.19'
$ws.Range('E37').Value = '  -5.26%  '
$ws.Range('E38').Value = '  +3.48%  '
Set-TextCell 'D39' '11.74'
$ws.Range('E39').Value = '  -5.81%  '
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('E41').Value = '  -1.23%  '
Set-TextCell 'D42' '59.11'
$ws.Range('E42').Value = '  -2.17%  '
$ws.Range('E43').Value = '  -6.98%  '
Set-TextCell 'D44' '101.40'
$ws.Range('E44').Value = '  +2.05%  '
Set-TextCell 'D45' '2.48'
$ws.Range('E45').Value = '  +9.18%  '
Set-TextCell 'D46' '0.0975'
$ws.Range('E46').Value = '  -1.55%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D47' '8.18'
$ws.Range('E47').Value = '  -4.57%  '
$ws.Range('B48').Value = 'WOONetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextCell 'D48' '0.458'
$ws.Range('E48').Value = '  +4.24%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('E51').Value = '  +0.14%  '
